$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(64, 8).Value = 7033.375
$ws.Cells.Item(64, 9).Value = 4073.625
$ws.Cells.Item(64, 10).Value = 9993.125
$ws.Cells.Item(64, 11).Value = 4073.625
$ws.Cells.Item(64, 12).Value = 9993.125
$ws.Cells.Item(64, 13).Value = -3825.625
$ws.Cells.Item(64, 14).Value = -10489.125
$ws.Cells.Item(67, 8).Value = 7033.375
$ws.Cells.Item(67, 9).Value = 4073.625
$ws.Cells.Item(67, 10).Value = 9993.125
$ws.Cells.Item(67, 11).Value = 4073.625
$ws.Cells.Item(67, 12).Value = 9993.125
$ws.Cells.Item(67, 13).Value = -3215.625
$ws.Cells.Item(67, 14).Value = -11709.125
$ws.Cells.Item(80, 8).Value = 909.7273
$ws.Cells.Item(80, 9).Value = 341.8
$ws.Cells.Item(80, 10).Value = 1383
$ws.Cells.Item(80, 11).Value = 1025.4
$ws.Cells.Item(80, 12).Value = 4149
$ws.Cells.Item(80, 13).Value = -27.40000000000009
$ws.Cells.Item(80, 14).Value = -6145
$ws.Cells.Item(83, 8).Value = 909.7273
$ws.Cells.Item(83, 9).Value = 341.8
$ws.Cells.Item(83, 10).Value = 1383
$ws.Cells.Item(83, 11).Value = 3076.2
$ws.Cells.Item(83, 12).Value = 12447
$ws.Cells.Item(83, 13).Value = 1915.8
$ws.Cells.Item(83, 14).Value = -22431
$ws.Cells.Item(88, 8).Value = 631.35297
$ws.Cells.Item(88, 10).Value = 727.75
$ws.Cells.Item(88, 12).Value = 727.75
$ws.Cells.Item(88, 14).Value = -1539.75
$ws.Cells.Item(91, 8).Value = 631.35297
$ws.Cells.Item(91, 10).Value = 727.75
$ws.Cells.Item(91, 12).Value = 727.75
$ws.Cells.Item(91, 14).Value = -3535.75
$ws.Cells.Item(129, 8).Value = 2292.2307
$ws.Cells.Item(129, 9).Value = 0
$ws.Cells.Item(129, 11).Value = 0
$ws.Cells.Item(129, 13).Value = $null
$ws.Cells.Item(137, 8).Value = 5688397.5
$ws.Cells.Item(137, 9).Value = 13161311
$ws.Cells.Item(137, 10).Value = 8983.68
$ws.Cells.Item(137, 11).Value = 39483933
$ws.Cells.Item(137, 12).Value = 26951.04
$ws.Cells.Item(137, 13).Value = -39481383
$ws.Cells.Item(137, 14).Value = -32051.04
$ws.Cells.Item(138, 8).Value = 6885.7144
$ws.Cells.Item(138, 9).Value = 3500
$ws.Cells.Item(138, 10).Value = 7450
$ws.Cells.Item(138, 11).Value = 10500
$ws.Cells.Item(138, 12).Value = 22350
$ws.Cells.Item(138, 13).Value = -5360
$ws.Cells.Item(138, 14).Value = -32630
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(45, 8).Value = 1089.0714
$ws.Cells.Item(45, 9).Value = 1020.5833
$ws.Cells.Item(45, 11).Value = 1020.5833
$ws.Cells.Item(45, 13).Value = -643.5833
$ws.Cells.Item(61, 8).Value = 5995.879
$ws.Cells.Item(61, 9).Value = 3698.8518
$ws.Cells.Item(61, 11).Value = 3698.8518
$ws.Cells.Item(61, 13).Value = -3486.8518
$ws.Cells.Item(131, 8).Value = 72500
$ws.Cells.Item(131, 10).Value = 72500
$ws.Cells.Item(131, 12).Value = 72500
$ws.Cells.Item(131, 14).Value = -82580
$ws.Cells.Item(136, 8).Value = 5995.879
$ws.Cells.Item(136, 9).Value = 3698.8518
$ws.Cells.Item(136, 11).Value = 11096.5554
$ws.Cells.Item(136, 13).Value = -8546.555399999999
$ws.Cells.Item(141, 8).Value = 119999.5
$ws.Cells.Item(141, 10).Value = 119999.5
$ws.Cells.Item(141, 12).Value = 119999.5
$ws.Cells.Item(141, 14).Value = -130359.5
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(97, 8).Value = 11397
$ws.Cells.Item(97, 9).Value = 7974.25
$ws.Cells.Item(97, 11).Value = 7974.25
$ws.Cells.Item(97, 13).Value = -6983.25
$ws.Cells.Item(107, 8).Value = 1484.909
$ws.Cells.Item(107, 10).Value = 1461.625
$ws.Cells.Item(107, 12).Value = 1461.625
$ws.Cells.Item(107, 14).Value = -5301.625
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(29, 8).Value = 6568.5293
$ws.Cells.Item(29, 9).Value = 4590
$ws.Cells.Item(29, 10).Value = 6692.1875
$ws.Cells.Item(29, 11).Value = 4590
$ws.Cells.Item(29, 12).Value = 6692.1875
$ws.Cells.Item(29, 13).Value = -4297
$ws.Cells.Item(29, 14).Value = -7278.1875
$ws.Cells.Item(70, 8).Value = 50000
$ws.Cells.Item(70, 10).Value = 50000
$ws.Cells.Item(70, 12).Value = 50000
$ws.Cells.Item(70, 14).Value = -50630
$ws.Cells.Item(73, 8).Value = 50000
$ws.Cells.Item(73, 10).Value = 50000
$ws.Cells.Item(73, 12).Value = 50000
$ws.Cells.Item(73, 14).Value = -52184
$ws.Cells.Item(107, 8).Value = 2411.647
$ws.Cells.Item(107, 10).Value = 3782.1667
$ws.Cells.Item(107, 12).Value = 3782.1667
$ws.Cells.Item(107, 14).Value = -7622.1667
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(134, 8).Value = 873.36365
$ws.Cells.Item(134, 9).Value = 873.36365
$ws.Cells.Item(134, 11).Value = 2620.09095
$ws.Cells.Item(134, 13).Value = 2449.90905
$ws.Cells.Item(136, 8).Value = 745
$ws.Cells.Item(136, 9).Value = 745
$ws.Cells.Item(136, 11).Value = 2235
$ws.Cells.Item(136, 13).Value = 2865
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(33, 8).Value = 1269498
$ws.Cells.Item(33, 10).Value = 25997.334
$ws.Cells.Item(33, 12).Value = 25997.334
$ws.Cells.Item(33, 14).Value = -26501.334
$ws.Cells.Item(36, 8).Value = 1516.5
$ws.Cells.Item(36, 9).Value = 599.2857
$ws.Cells.Item(36, 11).Value = 599.2857
$ws.Cells.Item(36, 13).Value = -114.2857
$ws.Cells.Item(43, 8).Value = 2150.25
$ws.Cells.Item(43, 9).Value = 2150.25
$ws.Cells.Item(43, 11).Value = 2150.25
$ws.Cells.Item(43, 13).Value = -1999.25
$ws.Cells.Item(102, 8).Value = 4366.3335
$ws.Cells.Item(102, 9).Value = 4366.3335
$ws.Cells.Item(102, 11).Value = 4366.3335
$ws.Cells.Item(102, 13).Value = -2744.3335
$ws.Cells.Item(122, 8).Value = 3739.0344
$ws.Cells.Item(122, 9).Value = 4011.28
$ws.Cells.Item(122, 11).Value = 12033.84
$ws.Cells.Item(122, 13).Value = -9583.84
$ws.Cells.Item(126, 8).Value = 3972.8
$ws.Cells.Item(126, 9).Value = 1487.5
$ws.Cells.Item(126, 11).Value = 4462.5
$ws.Cells.Item(126, 13).Value = -1992.5
$ws.Cells.Item(138, 8).Value = 99950
$ws.Cells.Item(138, 10).Value = 99950
$ws.Cells.Item(138, 12).Value = 99950
$ws.Cells.Item(138, 14).Value = -110230
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 27552.562
$ws.Cells.Item(7, 9).Value = 29649.334
$ws.Cells.Item(7, 11).Value = 29649.334
$ws.Cells.Item(7, 13).Value = -29537.334
$ws.Cells.Item(22, 8).Value = 2846.9546
$ws.Cells.Item(22, 9).Value = 1850.6666
$ws.Cells.Item(22, 10).Value = 4429.294
$ws.Cells.Item(22, 11).Value = 1850.6666
$ws.Cells.Item(22, 12).Value = 4429.294
$ws.Cells.Item(22, 13).Value = -1555.6666
$ws.Cells.Item(22, 14).Value = -5019.294
$ws.Cells.Item(27, 8).Value = 2846.9546
$ws.Cells.Item(27, 9).Value = 1850.6666
$ws.Cells.Item(27, 10).Value = 4429.294
$ws.Cells.Item(27, 11).Value = 1850.6666
$ws.Cells.Item(27, 12).Value = 4429.294
$ws.Cells.Item(27, 13).Value = -1743.6666
$ws.Cells.Item(27, 14).Value = -4643.294
$ws.Cells.Item(59, 8).Value = 75000
$ws.Cells.Item(59, 10).Value = 75000
$ws.Cells.Item(59, 12).Value = 75000
$ws.Cells.Item(59, 14).Value = -76308
$ws.Cells.Item(61, 8).Value = 4080
$ws.Cells.Item(61, 9).Value = 3625
$ws.Cells.Item(61, 10).Value = 5900
$ws.Cells.Item(61, 11).Value = 3625
$ws.Cells.Item(61, 12).Value = 5900
$ws.Cells.Item(61, 13).Value = -3423
$ws.Cells.Item(61, 14).Value = -6304
$ws.Cells.Item(113, 8).Value = 4080
$ws.Cells.Item(113, 9).Value = 3625
$ws.Cells.Item(113, 10).Value = 5900
$ws.Cells.Item(113, 11).Value = 3625
$ws.Cells.Item(113, 12).Value = 5900
$ws.Cells.Item(113, 13).Value = -1455
$ws.Cells.Item(113, 14).Value = -10240
$ws.Cells.Item(126, 8).Value = 27552.562
$ws.Cells.Item(126, 9).Value = 29649.334
$ws.Cells.Item(126, 11).Value = 88948.00199999999
$ws.Cells.Item(126, 13).Value = -86478.00199999999
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(122, 8).Value = 2831.1072
$ws.Cells.Item(122, 9).Value = 2388
$ws.Cells.Item(122, 11).Value = 7164
$ws.Cells.Item(122, 13).Value = -4714
$ws.Cells.Item(126, 8).Value = 54122.35
$ws.Cells.Item(126, 9).Value = 67490.44
$ws.Cells.Item(126, 10).Value = 650
$ws.Cells.Item(126, 11).Value = 202471.32
$ws.Cells.Item(126, 12).Value = 1950
$ws.Cells.Item(126, 13).Value = -200001.32
$ws.Cells.Item(126, 14).Value = -6890
